$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a Text number format to the cells we are about to populate, so that
# numeric-looking values (student id, semester count, weight) are stored as
# text rather than being auto-converted to numbers -- this mirrors the new
# "@" (text) cell style added to the workbook's cellXfs table.
$ws.Range("A1:E2").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "학번"
$ws.Range("B1").Value = "학생명"
$ws.Range("C1").Value = "학생 학기수"
$ws.Range("D1").Value = "전체 가중치"
$ws.Range("E1").Value = "결과"

# Data row
$ws.Range("A2").Value = "22200000"
$ws.Range("B2").Value = "김한동"
$ws.Range("C2").Value = "5"
$ws.Range("D2").Value = "30"
$ws.Range("E2").Value = "100만원"

$ws.Range("E12").Select() | Out-Null
